$d = $word.ActiveDocument

$d.Content.Find.Execute("93×50=", $true, $false, $false, $false, $false, $true, 1, $false, "14×39=", 2) | Out-Null
$d.Content.Find.Execute("34×12=", $true, $false, $false, $false, $false, $true, 1, $false, "85×56=", 2) | Out-Null
$d.Content.Find.Execute("77×74=", $true, $false, $false, $false, $false, $true, 1, $false, "21×74=", 2) | Out-Null
$d.Content.Find.Execute("94×66=", $true, $false, $false, $false, $false, $true, 1, $false, "99×92=", 2) | Out-Null
$d.Content.Find.Execute("29×77=", $true, $false, $false, $false, $false, $true, 1, $false, "44×81=", 2) | Out-Null
$d.Content.Find.Execute("98×70=", $true, $false, $false, $false, $false, $true, 1, $false, "40×60=", 2) | Out-Null
$d.Content.Find.Execute("38×62=", $true, $false, $false, $false, $false, $true, 1, $false, "59×58=", 2) | Out-Null
$d.Content.Find.Execute("59×16=", $true, $false, $false, $false, $false, $true, 1, $false, "29×28=", 2) | Out-Null
$d.Content.Find.Execute("70×65=", $true, $false, $false, $false, $false, $true, 1, $false, "87×44=", 2) | Out-Null
$d.Content.Find.Execute("41×71=", $true, $false, $false, $false, $false, $true, 1, $false, "84×43=", 2) | Out-Null
$d.Content.Find.Execute("11×46=", $true, $false, $false, $false, $false, $true, 1, $false, "47×42=", 2) | Out-Null
$d.Content.Find.Execute("26×33=", $true, $false, $false, $false, $false, $true, 1, $false, "77×88=", 2) | Out-Null
$d.Content.Find.Execute("30×24=", $true, $false, $false, $false, $false, $true, 1, $false, "43×24=", 2) | Out-Null
$d.Content.Find.Execute("64×17=", $true, $false, $false, $false, $false, $true, 1, $false, "19×22=", 2) | Out-Null
$d.Content.Find.Execute("75×60=", $true, $false, $false, $false, $false, $true, 1, $false, "66×82=", 2) | Out-Null
$d.Content.Find.Execute("99×38=", $true, $false, $false, $false, $false, $true, 1, $false, "21×69=", 2) | Out-Null
$d.Content.Find.Execute("16×75=", $true, $false, $false, $false, $false, $true, 1, $false, "86×70=", 2) | Out-Null
$d.Content.Find.Execute("44×31=", $true, $false, $false, $false, $false, $true, 1, $false, "15×82=", 2) | Out-Null
$d.Content.Find.Execute("37×17=", $true, $false, $false, $false, $false, $true, 1, $false, "15×81=", 2) | Out-Null
$d.Content.Find.Execute("24×28=", $true, $false, $false, $false, $false, $true, 1, $false, "86×94=", 2) | Out-Null
$d.Content.Find.Execute("87×20=", $true, $false, $false, $false, $false, $true, 1, $false, "11×85=", 2) | Out-Null
$d.Content.Find.Execute("47×97=", $true, $false, $false, $false, $false, $true, 1, $false, "35×13=", 2) | Out-Null
$d.Content.Find.Execute("36×27=", $true, $false, $false, $false, $false, $true, 1, $false, "66×42=", 2) | Out-Null
$d.Content.Find.Execute("58×71=", $true, $false, $false, $false, $false, $true, 1, $false, "66×16=", 2) | Out-Null
$d.Content.Find.Execute("33×46=", $true, $false, $false, $false, $false, $true, 1, $false, "56×91=", 2) | Out-Null

Write-Host "Done applying replacements."
